# Scheduled-runner refresh of market-price-derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
# Values below are the freshly recomputed figures for the affected leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 931.6667
$ws.Range("I15").Value = 931.6667
$ws.Range("K15").Value = 2795.0001
$ws.Range("M15").Value = -2626.0001
$ws.Range("H39").Value = 97.38461
$ws.Range("I39").Value = 97.38461
$ws.Range("K39").Value = 292.15383
$ws.Range("M39").Value = 3.846170000000029
$ws.Range("H40").Value = 1780
$ws.Range("I40").Value = 1847.5
$ws.Range("K40").Value = 1847.5
$ws.Range("M40").Value = -1672.5
$ws.Range("H42").Value = 83.78570999999999
$ws.Range("I42").Value = 60.75
$ws.Range("K42").Value = 182.25
$ws.Range("M42").Value = 47.75
$ws.Range("H62").Value = 1974.3077
$ws.Range("I62").Value = 2531.75
$ws.Range("J62").Value = 1082.4
$ws.Range("K62").Value = 2531.75
$ws.Range("L62").Value = 1082.4
$ws.Range("M62").Value = -1907.75
$ws.Range("N62").Value = -2330.4
$ws.Range("H64").Value = 3417.2354
$ws.Range("I64").Value = 3341.4285
$ws.Range("J64").Value = 3470.3
$ws.Range("K64").Value = 3341.4285
$ws.Range("L64").Value = 3470.3
$ws.Range("M64").Value = -3093.4285
$ws.Range("N64").Value = -3966.3
$ws.Range("H65").Value = 1974.3077
$ws.Range("I65").Value = 2531.75
$ws.Range("J65").Value = 1082.4
$ws.Range("K65").Value = 12658.75
$ws.Range("L65").Value = 5412
$ws.Range("M65").Value = -9538.75
$ws.Range("N65").Value = -11652
$ws.Range("H67").Value = 3417.2354
$ws.Range("I67").Value = 3341.4285
$ws.Range("J67").Value = 3470.3
$ws.Range("K67").Value = 3341.4285
$ws.Range("L67").Value = 3470.3
$ws.Range("M67").Value = -2483.4285
$ws.Range("N67").Value = -5186.3
$ws.Range("H98").Value = 2079
$ws.Range("I98").Value = 2300.8333
$ws.Range("J98").Value = 1635.3334
$ws.Range("K98").Value = 2300.8333
$ws.Range("L98").Value = 1635.3334
$ws.Range("M98").Value = -802.8332999999998
$ws.Range("N98").Value = -4631.3334
$ws.Range("H122").Value = 2079
$ws.Range("I122").Value = 2300.8333
$ws.Range("J122").Value = 1635.3334
$ws.Range("K122").Value = 6902.499899999999
$ws.Range("L122").Value = 4906.0002
$ws.Range("M122").Value = -4452.499899999999
$ws.Range("N122").Value = -9806.0002
$ws.Range("H138").Value = 2261.7844
$ws.Range("I138").Value = 1430.8077
$ws.Range("J138").Value = 3126
$ws.Range("K138").Value = 4292.4231
$ws.Range("L138").Value = 9378
$ws.Range("M138").Value = 847.5769
$ws.Range("N138").Value = -19658
$ws.Range("H139").Value = 50172.25
$ws.Range("J139").Value = 50172.25
$ws.Range("L139").Value = 50172.25
$ws.Range("N139").Value = -60452.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15626504
$ws.Range("I61").Value = 20001606
$ws.Range("J61").Value = 1143.7142
$ws.Range("K61").Value = 20001606
$ws.Range("L61").Value = 1143.7142
$ws.Range("M61").Value = -20001394
$ws.Range("N61").Value = -1567.7142
$ws.Range("H74").Value = 14288954
$ws.Range("I74").Value = 20835044
$ws.Range("K74").Value = 20835044
$ws.Range("M74").Value = -20834170
$ws.Range("H77").Value = 14288954
$ws.Range("I77").Value = 20835044
$ws.Range("K77").Value = 104175220
$ws.Range("M77").Value = -104170852
$ws.Range("H136").Value = 15626504
$ws.Range("I136").Value = 20001606
$ws.Range("J136").Value = 1143.7142
$ws.Range("K136").Value = 60004818
$ws.Range("L136").Value = 3431.1426
$ws.Range("M136").Value = -60002268
$ws.Range("N136").Value = -8531.142599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5955628
$ws.Range("I31").Value = 3399.9058
$ws.Range("K31").Value = 3399.9058
$ws.Range("M31").Value = -3104.9058
$ws.Range("H34").Value = 5955628
$ws.Range("I34").Value = 3399.9058
$ws.Range("K34").Value = 3399.9058
$ws.Range("M34").Value = -3197.9058
$ws.Range("H105").Value = 1160
$ws.Range("I105").Value = 990
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 990
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 757
$ws.Range("N105").Value = -4994
$ws.Range("H132").Value = 10871156
$ws.Range("I132").Value = 13890132
$ws.Range("J132").Value = 2840.9
$ws.Range("K132").Value = 41670396
$ws.Range("L132").Value = 8522.700000000001
$ws.Range("M132").Value = -41667866
$ws.Range("N132").Value = -13582.7
$ws.Range("H140").Value = 35933.562
$ws.Range("J140").Value = 35933.562
$ws.Range("L140").Value = 35933.562
$ws.Range("N140").Value = -46293.562
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 823.1799999999999
$ws.Range("I131").Value = 467.77777
$ws.Range("J131").Value = 858.32965
$ws.Range("K131").Value = 1403.33331
$ws.Range("L131").Value = 2574.98895
$ws.Range("M131").Value = 3636.66669
$ws.Range("N131").Value = -12654.98895
$ws.Range("H137").Value = 4722.355
$ws.Range("I137").Value = 2887.6428
$ws.Range("J137").Value = 6233.294
$ws.Range("K137").Value = 8662.928400000001
$ws.Range("L137").Value = 18699.882
$ws.Range("M137").Value = -3562.928400000001
$ws.Range("N137").Value = -28899.882
$ws.Range("H140").Value = 3043.889
$ws.Range("I140").Value = 1652.6666
$ws.Range("J140").Value = 10000
$ws.Range("K140").Value = 4957.9998
$ws.Range("L140").Value = 30000
$ws.Range("M140").Value = 222.0002000000004
$ws.Range("N140").Value = -40360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 66666664
$ws.Range("I122").Value = 66666664
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 199999992
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -199997542
$ws.Range("N122").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11875
$ws.Range("I40").Value = 18666.666
$ws.Range("J40").Value = 7800
$ws.Range("K40").Value = 18666.666
$ws.Range("L40").Value = 7800
$ws.Range("M40").Value = -18530.666
$ws.Range("N40").Value = -8072
$ws.Range("H122").Value = 10637
$ws.Range("I122").Value = 27200
$ws.Range("J122").Value = 6496.25
$ws.Range("K122").Value = 81600
$ws.Range("L122").Value = 19488.75
$ws.Range("M122").Value = -79150
$ws.Range("N122").Value = -24388.75
$ws.Range("H132").Value = 7581484
$ws.Range("I132").Value = 3347.2046
$ws.Range("J132").Value = 22737758
$ws.Range("K132").Value = 10041.6138
$ws.Range("L132").Value = 68213274
$ws.Range("M132").Value = -7511.613799999999
$ws.Range("N132").Value = -68218334
$ws.Range("H136").Value = 11368120
$ws.Range("I136").Value = 15153172
$ws.Range("J136").Value = 12964.546
$ws.Range("K136").Value = 45459516
$ws.Range("L136").Value = 38893.638
$ws.Range("M136").Value = -45456966
$ws.Range("N136").Value = -43993.638
$ws.Range("H139").Value = 58682.668
$ws.Range("J139").Value = 58682.668
$ws.Range("L139").Value = 58682.668
$ws.Range("N139").Value = -68962.66800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1262.1936
$ws.Range("I136").Value = 1114
$ws.Range("K136").Value = 3342
$ws.Range("M136").Value = -792